$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 2376
$ws1.Range("F11").Value = 560
$ws1.Range("F15").Value = 188
$ws1.Range("F21").Value = 2596
$ws1.Range("F23").Value = 22
$ws1.Range("F27").Value = 1735
$ws1.Range("F33").Value = 4512

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F24").Value = 193

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F16").Value = 2376
$ws4.Range("F18").Value = 560
$ws4.Range("F23").Value = 188
$ws4.Range("F32").Value = 2596
$ws4.Range("F35").Value = 22
$ws4.Range("F41").Value = 1735
$ws4.Range("F47").Value = 4512
